$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A (Serie) cells for the new rows to be written as plain text
# (not auto-converted to date serials) by pre-setting NumberFormat to "@",
# then restoring the default "Normal" style afterwards so no residual
# number-format styling is left behind on the cells.
$dateRange = $ws.Range("A198:A203")
$dateRange.NumberFormat = "@"

# Row 197
$ws.Cells.Item(197, 2).Value = 34002.9
$ws.Cells.Item(197, 3).Value = 14255.5
$ws.Cells.Item(197, 4).Value = 7011
$ws.Cells.Item(197, 5).Value = 28444.9
$ws.Cells.Item(197, 6).Value = 6477.7
$ws.Cells.Item(197, 7).Value = 15036.6
$ws.Cells.Item(197, 10).Value = 1522.5
$ws.Cells.Item(197, 11).Value = 71265.4
$ws.Cells.Item(197, 12).Value = 4112.7
$ws.Cells.Item(197, 13).Value = 1614.5
$ws.Cells.Item(197, 14).Value = 16408.4
$ws.Cells.Item(197, 15).Value = 1394.4
$ws.Cells.Item(197, 16).Value = 77578.9
$ws.Cells.Item(197, 17).Value = 110393.1
$ws.Cells.Item(197, 18).Value = 1381.8
$ws.Cells.Item(197, 19).Value = 50862.9
$ws.Cells.Item(197, 20).Value = 18192.1

# Row 198
$ws.Cells.Item(198, 1).Value = "05-10-2021"
$ws.Cells.Item(198, 2).Value = 34314.7
$ws.Cells.Item(198, 3).Value = 14433.8
$ws.Cells.Item(198, 4).Value = 7077.1
$ws.Cells.Item(198, 5).Value = 27822.1
$ws.Cells.Item(198, 6).Value = 6576.3
$ws.Cells.Item(198, 7).Value = 15194.5
$ws.Cells.Item(198, 8).Value = 2962.2
$ws.Cells.Item(198, 10).Value = 1530.4
$ws.Cells.Item(198, 11).Value = 71940.5
$ws.Cells.Item(198, 12).Value = 4223.9
$ws.Cells.Item(198, 13).Value = 1624.2
$ws.Cells.Item(198, 14).Value = 16460.8
$ws.Cells.Item(198, 15).Value = 1391.7
$ws.Cells.Item(198, 16).Value = 78191.5
$ws.Cells.Item(198, 17).Value = 110457.6
$ws.Cells.Item(198, 18).Value = 1383.2
$ws.Cells.Item(198, 19).Value = 51056.9
$ws.Cells.Item(198, 20).Value = 18342.7

# Row 199
$ws.Cells.Item(199, 1).Value = "06-10-2021"
$ws.Cells.Item(199, 2).Value = 34417
$ws.Cells.Item(199, 3).Value = 14501.9
$ws.Cells.Item(199, 4).Value = 6995.9
$ws.Cells.Item(199, 5).Value = 27528.9
$ws.Cells.Item(199, 6).Value = 6493.1
$ws.Cells.Item(199, 7).Value = 14973.3
$ws.Cells.Item(199, 8).Value = 2908.3
$ws.Cells.Item(199, 10).Value = 1559.4
$ws.Cells.Item(199, 11).Value = 71789.3
$ws.Cells.Item(199, 12).Value = 4171.7
$ws.Cells.Item(199, 13).Value = 1619.5
$ws.Cells.Item(199, 14).Value = 16393.2
$ws.Cells.Item(199, 15).Value = 1371.7
$ws.Cells.Item(199, 16).Value = 77443.1
$ws.Cells.Item(199, 17).Value = 110559.6
$ws.Cells.Item(199, 18).Value = 1381.9
$ws.Cells.Item(199, 19).Value = 50974.5
$ws.Cells.Item(199, 20).Value = 18356.1

# Row 200
$ws.Cells.Item(200, 1).Value = "07-10-2021"
$ws.Cells.Item(200, 2).Value = 34754.9
$ws.Cells.Item(200, 3).Value = 14654
$ws.Cells.Item(200, 4).Value = 7078
$ws.Cells.Item(200, 5).Value = 27678.2
$ws.Cells.Item(200, 6).Value = 6600.2
$ws.Cells.Item(200, 7).Value = 15250.9
$ws.Cells.Item(200, 8).Value = 2959.5
$ws.Cells.Item(200, 10).Value = 1561.3
$ws.Cells.Item(200, 11).Value = 72124.9
$ws.Cells.Item(200, 12).Value = 4228
$ws.Cells.Item(200, 13).Value = 1633.7
$ws.Cells.Item(200, 14).Value = 16713.9
$ws.Cells.Item(200, 15).Value = 1396.9
$ws.Cells.Item(200, 16).Value = 77553.3
$ws.Cells.Item(200, 17).Value = 110585.4
$ws.Cells.Item(200, 18).Value = 1393.5
$ws.Cells.Item(200, 19).Value = 50877
$ws.Cells.Item(200, 20).Value = 19375.7

# Row 201
$ws.Cells.Item(201, 1).Value = "08-10-2021"
$ws.Cells.Item(201, 2).Value = 34746.3
$ws.Cells.Item(201, 3).Value = 14579.5
$ws.Cells.Item(201, 4).Value = 7095.6
$ws.Cells.Item(201, 5).Value = 28048.9
$ws.Cells.Item(201, 6).Value = 6560
$ws.Cells.Item(201, 7).Value = 15206.1
$ws.Cells.Item(201, 8).Value = 2956.3
$ws.Cells.Item(201, 9).Value = 4929.9
$ws.Cells.Item(201, 10).Value = 1563.9
$ws.Cells.Item(201, 11).Value = 73327.7
$ws.Cells.Item(201, 12).Value = 4238
$ws.Cells.Item(201, 13).Value = 1639.4
$ws.Cells.Item(201, 14).Value = 16640.4
$ws.Cells.Item(201, 15).Value = 1398
$ws.Cells.Item(201, 17).Value = 112833.2
$ws.Cells.Item(201, 18).Value = 1401
$ws.Cells.Item(201, 19).Value = 51136.6

# Row 202
$ws.Cells.Item(202, 1).Value = "11-10-2021"
$ws.Cells.Item(202, 2).Value = 34496.1
$ws.Cells.Item(202, 3).Value = 14486.2
$ws.Cells.Item(202, 4).Value = 7146.9
$ws.Cells.Item(202, 5).Value = 28498.2
$ws.Cells.Item(202, 6).Value = 6570.5
$ws.Cells.Item(202, 7).Value = 15199.1
$ws.Cells.Item(202, 9).Value = 4936.2
$ws.Cells.Item(202, 10).Value = 1570.8
$ws.Cells.Item(202, 11).Value = 74459.9
$ws.Cells.Item(202, 12).Value = 4279.9
$ws.Cells.Item(202, 13).Value = 1633.4
$ws.Cells.Item(202, 15).Value = 1416.3
$ws.Cells.Item(202, 17).Value = 112180.5
$ws.Cells.Item(202, 18).Value = 1406.3
$ws.Cells.Item(202, 19).Value = 51646.8
$ws.Cells.Item(202, 20).Value = 19585.6

# Row 203
$ws.Cells.Item(203, 1).Value = "12-10-2021"
$ws.Cells.Item(203, 4).Value = 7113.3
$ws.Cells.Item(203, 5).Value = 28230.6
$ws.Cells.Item(203, 6).Value = 6536.1
$ws.Cells.Item(203, 7).Value = 15146.5
$ws.Cells.Item(203, 8).Value = 2916.4
$ws.Cells.Item(203, 9).Value = 4883.8
$ws.Cells.Item(203, 10).Value = 1583.9
$ws.Cells.Item(203, 11).Value = 74216.1
$ws.Cells.Item(203, 12).Value = 4266
$ws.Cells.Item(203, 13).Value = 1643.6
$ws.Cells.Item(203, 14).Value = 16462.8
$ws.Cells.Item(203, 15).Value = 1420.2

# Restore default styling on the date column cells (removes the temporary
# text number-format applied above).
$dateRange.Style = "Normal"

